$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.050.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.83%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.016.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.99%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'226.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.80%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.600"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.69%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'54.68"
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = "'  -3.90%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0778"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.40%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.102"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -5.46%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.315.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.92%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D14").Value = "'20.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -5.11%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.740"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.37%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'5.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.26%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.020.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.71%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'37.008.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'6.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.73%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'68.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0817"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.66%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'223.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.05%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "'  +1.91%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -7.48%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'166.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.03%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -7.68%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -2.48%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'18.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.40%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -6.64%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -3.65%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.52%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -3.03%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -5.30%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -7.43%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.95%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.09%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -5.40%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.03%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.473.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.18%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -5.30%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'94.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.73%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -4.98%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'16.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.08%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -5.20%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -6.19%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'ARBITRUM"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'1.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.36%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'FraxShare"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'7.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.32%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -1.87%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.201.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.03%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'FTXToken"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'3.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -13.64%  "
$ws.Range("E51").Style = "Normal"
